# Apply crypto price/volume refresh (and a few row re-orderings) per commit
# "Updated cryptos list on Sun Oct  8 11:58:18 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.829.06'
$ws.Range('E2').Value = '  -0.49%  '

$ws.Range('D3').Value = '1.617.56'
$ws.Range('E3').Value = '  -1.36%  '

$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.09%  '

$ws.Range('D5').Value = "'210.02"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.30%  '

$ws.Range('E7').Value = '  -0.11%  '

$ws.Range('D8').Value = "'23.22"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.01%  '

$ws.Range('E9').Value = '  -0.77%  '

$ws.Range('E10').Value = '  -0.80%  '

$ws.Range('D12').Value = '1.847.11'
$ws.Range('E12').Value = '  -1.38%  '

$ws.Range('D13').Value = '1.612.50'
$ws.Range('E13').Value = '  -1.67%  '

$ws.Range('E14').Value = '  -1.94%  '

$ws.Range('E15').Value = '  -2.87%  '

$ws.Range('D16').Value = "'65.08"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.07%  '

$ws.Range('D17').Value = '27.814.29'
$ws.Range('E17').Value = '  -0.53%  '

$ws.Range('E18').Value = '  -2.74%  '

$ws.Range('D19').Value = '0.0₃0718'
$ws.Range('E19').Value = '  -0.77%  '

$ws.Range('D20').Value = "'7.58"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.21%  '

$ws.Range('E21').Value = '  -0.15%  '

$ws.Range('E22').Value = '  -1.72%  '

$ws.Range('D23').Value = "'10.04"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.11%  '

$ws.Range('E24').Value = '  -3.24%  '

$ws.Range('D25').Value = "'154.28"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.63%  '

$ws.Range('E26').Value = '  -1.24%  '

$ws.Range('E27').Value = '  -0.17%  '

$ws.Range('B28').Value = 'BinanceUSD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D28').Value = "'1.00"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.14%  '

$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = "'15.43"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.43%  '

$ws.Range('E30').Value = '  -1.51%  '

$ws.Range('E31').Value = '  -0.92%  '

$ws.Range('D32').Value = "'3.43"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.25%  '

$ws.Range('E33').Value = '  -1.65%  '

$ws.Range('D34').Value = '1.384.76'
$ws.Range('E34').Value = '  -2.08%  '

$ws.Range('E35').Value = '  -0.51%  '

$ws.Range('D36').Value = "'0.996"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +10.53%  '

$ws.Range('E37').Value = '  -1.10%  '

$ws.Range('E38').Value = '  +0.02%  '

$ws.Range('D39').Value = "'0.553"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.83%  '

$ws.Range('D40').Value = "'0.841"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.61%  '

$ws.Range('E41').Value = '  -0.17%  '

$ws.Range('D42').Value = "'0.991"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.41%  '

$ws.Range('E43').Value = '  -1.55%  '

$ws.Range('D44').Value = "'5.45"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.34%  '

$ws.Range('D45').Value = "'65.39"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.78%  '

$ws.Range('D46').Value = '1.757.46'
$ws.Range('E46').Value = '  -1.42%  '

$ws.Range('E47').Value = '  -2.26%  '

$ws.Range('D48').Value = "'87.54"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.20%  '

$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0103'
$ws.Range('E49').Value = '  -2.34%  '

$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = "'0.101"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.46%  '

$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = "'0.0502"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.73%  '
